# Lab01_ReviewReport.xlsx -- "Review report + corrections"
#
# Fills in the three per-phase review-defect sheets with the actual
# reviewer, date, and defect rows, sets the "effort to review" duration,
# removes the stray spacer cell above the effort row, and restores the
# selection/active-sheet state recorded in the authored workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Requirements Phase Defects"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Requirements Phase Defects")

$ws1.Range("D6").Value = "Bulmez Alexandru - Florin"
$ws1.Range("D7").Value = "03.18.2019"

$ws1.Rows.Item(10).RowHeight = 30
$ws1.Range("C10").Value = "R1,R2"
$ws1.Range("D10").Value = "par 1"
$ws1.Range("E10").Value = " Nu se precizeaza tipul de aplicatie  (aplicatie consola sau aplicatie cu interfata grafica)"

$ws1.Rows.Item(11).RowHeight = 30
$ws1.Range("C11").Value = "R1"
$ws1.Range("D11").Value = "par 1"
$ws1.Range("E11").Value = "Nu se precizeaza formatul si numele fisierelor de intrare"

$ws1.Rows.Item(12).RowHeight = 30
$ws1.Range("C12").Value = "R1"
$ws1.Range("D12").Value = "par 5"
$ws1.Range("E12").Value = "Nu se precizeaza conditia conform careia un elev este considerat corigent la o materie"

# Stray empty spacer cell above the "effort" row goes away entirely.
$ws1.Range("E26").Clear()

# Effort to review document: 30 minutes.
$ws1.Range("E27").NumberFormat = "h:mm"
$ws1.Range("E27").Value = 0.020833333333333332

$ws1.Range("E11").Select()

# ---------------------------------------------------------------------
# Sheet 2: "Architect. Design Phase Defects"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Architect. Design Phase Defects")

$ws2.Range("D6").Value = "Bulmez Alexandru - Florin"
$ws2.Range("D7").Value = "03.18.2019"

$ws2.Rows.Item(10).RowHeight = 30
$ws2.Range("C10").Value = "A10"
$ws2.Range("E10").Value = "Nu, lipsesc anumite concepte: conceptul de Nota"

$ws2.Rows.Item(11).RowHeight = 45
$ws2.Range("C11").Value = "A5"
$ws2.Range("E11").Value = "Nu, desi prezenta in diagrama, clasa ClasaException nu este utilizata (conform relatiilor dintre clase)"

$ws2.Rows.Item(12).RowHeight = 30
$ws2.Range("C12").Value = "A2"
$ws2.Range("E12").Value = "Nu, lipseste stratul de UI din cadrul unei arhitecturi stratificate"

# Stray empty spacer cell above the "effort" row goes away entirely.
$ws2.Range("E27").Clear()

# Effort to review document: 30 minutes.
$ws2.Range("E28").NumberFormat = "h:mm"
$ws2.Range("E28").Value = 0.020833333333333332

$ws2.Activate()
$ws2.Range("E13").Select()

# ---------------------------------------------------------------------
# Sheet 3: "Coding Phase Defects"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Coding Phase Defects")

$ws3.Range("D6").Value = "Bulmez Alexandru - Florin"
$ws3.Range("D7").Value = "03.18.2019"

$ws3.Rows.Item(10).RowHeight = 60
$ws3.Range("C10").Value = "C1"
$ws3.Range("E10").Value = "Clasa ClasaException este importata gresit in cadrul claselor StartApp, AddNotaTest, CalculeazaMediiTest, GetCorigentiTest, IntegrationTest, NoteController"

$ws3.Range("C11").Value = "C4"
$ws3.Range("E11").Value = "Lipsa fisierelor de intrare pentru elevi si note"

$ws3.Rows.Item(12).RowHeight = 45
$ws3.Range("C12").Value = "C1"
$ws3.Range("E12").Value = "In cadrul clasei Elev campul nrmatricol este de tip int, pe cand in cadrul clasei Nota campul nrmatricol este de tip double"

# Stray empty spacer cell above the "effort" row goes away entirely.
$ws3.Range("E31").Clear()

# Effort to review document: 1 hour.
$ws3.Range("E32").NumberFormat = "h:mm"
$ws3.Range("E32").Value = 0.041666666666666664

$ws3.Range("C37").Select()

# ---------------------------------------------------------------------
# The last-active sheet when the workbook was saved was sheet 2.
# ---------------------------------------------------------------------
$ws2.Activate()
